# Updates the cryptocurrency price/volume table with the latest pulled values.
# (GitHub Actions data-refresh commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference a cell whose style is (and always stays) the default/unstyled
# cell style, so numeric-looking text values below can be re-stamped back
# onto the plain style after Excel auto-applies a quote-prefix style to them.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "62.594.81"
$ws.Range("E2").Value = "  +5.35%  "
$ws.Range("D3").Value = "3.468.65"
$ws.Range("E3").Value = "  +4.57%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'408.23"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'129.16"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +15.87%  "
$ws.Range("D7").Value = "3.461.73"
$ws.Range("E7").Value = "  +4.63%  "
$ws.Range("D8").Value = "'0.595"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.691"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +9.71%  "
$ws.Range("D11").Value = "'0.126"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +28.41%  "
$ws.Range("D12").Value = "'42.66"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +7.61%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "4.014.84"
$ws.Range("E14").Value = "  +4.67%  "
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "'20.00"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").Value = "3.511.45"
$ws.Range("E17").Value = "  +6.03%  "
$ws.Range("D18").Value = "62.541.01"
$ws.Range("E18").Value = "  +5.66%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").Value = "'0.0000135"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +22.34%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'82.35"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +9.35%  "
$ws.Range("D24").Value = "'13.10"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'309.64"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "'30.30"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +6.47%  "
$ws.Range("D28").Value = "'8.23"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +5.40%  "
$ws.Range("E29").Value = "  +5.32%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "'4.37"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("D33").Value = "'11.91"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("D35").Value = "'43.28"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  +9.08%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "'0.0493"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("D38").Value = "'52.56"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'2.96"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -5.23%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.99"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +3.90%  "
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'137.59"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "'17.50"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +4.24%  "
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "'3.96"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D48").Value = "'2.26"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").Value = "'22.45"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "2.210.51"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "3.815.90"
$ws.Range("E51").Value = "  +4.87%  "
